# Generate Report for Handback
#
# This reproduces the "handback" localization-status update:
#  - Overview / zh-cn / de-de "Status" cells move from "Ready for handoff"
#    to "Handed back: in sync with en-US"
#  - zh-cn gains its "Latest Target File" (hyperlinked .md) + "Latest
#    Handback File" (.xlf) values on both data rows, and its existing
#    "Latest Handback DateTime" placeholder is replaced with a real
#    timestamp.
#  - de-de gets the same "Latest Target File" / "Latest Handback File"
#    treatment, plus a (later) handback timestamp of its own.
#  - A few columns that now hold longer filenames are widened to fit.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdFile1 = "78736c29-08b4-4075-865f-7cd8f4a890a7.md"
$mdFile2 = "bf1b051a-e01c-41fe-b408-9333c50ee89f.md"
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9601d031f74a590d14594e2214385e17a19672bd/e2e/78736c29-08b4-4075-865f-7cd8f4a890a7.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9601d031f74a590d14594e2214385e17a19672bd/e2e/bf1b051a-e01c-41fe-b408-9333c50ee89f.md"

# ---------------------------------------------------------------
# Overview sheet: just the "Status" text + widened zh-cn/de-de cols
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("I2").Value = $mdFile1
$wsZh.Range("J2").Value = "78736c29-08b4-4075-865f-7cd8f4a890a7.317daa7e477ee8495b0cd4769b5d3ce58961cbad.zh-cn.xlf"

$wsZh.Range("I3").Value = $mdFile2
$wsZh.Range("J3").Value = "bf1b051a-e01c-41fe-b408-9333c50ee89f.a40316f5666f25638bb2c79f8598128643f7bfc5.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-08-18 02:24:53"
$wsZh.Range("K3").Value = "2016-08-18 02:24:53"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFile1)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl2, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFile2)

$wsZh.Range("I2").Style = "Hyperlink"
$wsZh.Range("I3").Style = "Hyperlink"

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("I2").Value = $mdFile1
$wsDe.Range("J2").Value = "78736c29-08b4-4075-865f-7cd8f4a890a7.317daa7e477ee8495b0cd4769b5d3ce58961cbad.de-de.xlf"

$wsDe.Range("I3").Value = $mdFile2
$wsDe.Range("J3").Value = "bf1b051a-e01c-41fe-b408-9333c50ee89f.a40316f5666f25638bb2c79f8598128643f7bfc5.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-18 02:25:00"
$wsDe.Range("K3").Value = "2016-08-18 02:25:00"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFile1)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl2, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFile2)

$wsDe.Range("I2").Style = "Hyperlink"
$wsDe.Range("I3").Style = "Hyperlink"

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664
